$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = -0.65
$ws.Range("G3").Value = -0.64
$ws.Range("E5").Value = 0.01
$ws.Range("F5").Value = 0.001
$ws.Range("G6").Value = -0.61
$ws.Range("G7").Value = -0.59
$ws.Range("G8").Value = -0.57
$ws.Range("G9").Value = -0.55
$ws.Range("I9").Value = 0
$ws.Range("G10").Value = -0.53
$ws.Range("G11").Value = -0.51
$ws.Range("I12").Value = 0
$ws.Range("E13").Value = 0.03
$ws.Range("F13").Value = 0.0021
$ws.Range("G14").Value = -0.46
$ws.Range("E15").Value = 0.04
$ws.Range("F15").Value = 0.0027
$ws.Range("G15").Value = -0.44
$ws.Range("H15").Value = 0.01
$ws.Range("G16").Value = -0.43
$ws.Range("H16").Value = 0.02
$ws.Range("G17").Value = -0.42
$ws.Range("H17").Value = 0.03
$ws.Range("G18").Value = -0.41
$ws.Range("H18").Value = 0.04
$ws.Range("E19").Value = 0.05
$ws.Range("F19").Value = 0.0031
$ws.Range("G19").Value = -0.4
$ws.Range("H19").Value = 0.06
$ws.Range("I19").Value = 0
$ws.Range("G20").Value = -0.39
$ws.Range("H20").Value = 0.09
$ws.Range("E21").Value = 0.06
$ws.Range("F21").Value = 0.0036
$ws.Range("G21").Value = -0.38
$ws.Range("H21").Value = 0.12
$ws.Range("G22").Value = -0.37
$ws.Range("H22").Value = 0.17
$ws.Range("E23").Value = 0.06
$ws.Range("F23").Value = 0.0035
$ws.Range("G23").Value = -0.36
$ws.Range("H23").Value = 0.23
$ws.Range("E24").Value = 0.07000000000000001
$ws.Range("F24").Value = 0.0041
$ws.Range("G24").Value = -0.35
$ws.Range("H24").Value = 0.3
$ws.Range("E25").Value = 0.08
$ws.Range("F25").Value = 0.0046
$ws.Range("G25").Value = -0.34
$ws.Range("H25").Value = 0.4
$ws.Range("E26").Value = 0.08
$ws.Range("F26").Value = 0.0044
$ws.Range("G26").Value = -0.32
$ws.Range("H26").Value = 0.68
$ws.Range("E27").Value = 0.09
$ws.Range("F27").Value = 0.0049
$ws.Range("G27").Value = -0.31
$ws.Range("H27").Value = 0.87
$ws.Range("E28").Value = 0.09
$ws.Range("F28").Value = 0.0049
$ws.Range("G28").Value = -0.31
$ws.Range("H28").Value = 1.09
$ws.Range("G29").Value = -0.3
$ws.Range("H29").Value = 1.37
$ws.Range("E30").Value = 0.12
$ws.Range("F30").Value = 0.0063
$ws.Range("G30").Value = -0.29
$ws.Range("H30").Value = 1.7
$ws.Range("I30").Value = 0
$ws.Range("G31").Value = -0.28
$ws.Range("H31").Value = 2.09
$ws.Range("I31").Value = 0
$ws.Range("E32").Value = 0.13
$ws.Range("F32").Value = 0.0067
$ws.Range("G32").Value = -0.27
$ws.Range("H32").Value = 2.54
$ws.Range("E33").Value = 0.13
$ws.Range("F33").Value = 0.0066
$ws.Range("G33").Value = -0.26
$ws.Range("H33").Value = 3.07
$ws.Range("I33").Value = 0
$ws.Range("E34").Value = 0.16
$ws.Range("F34").Value = 0.008
$ws.Range("G34").Value = -0.25
$ws.Range("H34").Value = 3.68
$ws.Range("E35").Value = 0.17
$ws.Range("F35").Value = 0.008399999999999999
$ws.Range("G35").Value = -0.24
$ws.Range("H35").Value = 4.37
$ws.Range("I35").Value = 0
$ws.Range("E36").Value = 0.17
$ws.Range("F36").Value = 0.0083
$ws.Range("G36").Value = -0.23
$ws.Range("H36").Value = 5.16
$ws.Range("I36").Value = 0
$ws.Range("E37").Value = 0.17
$ws.Range("F37").Value = 0.008200000000000001
$ws.Range("G37").Value = -0.22
$ws.Range("H37").Value = 6.04
$ws.Range("E38").Value = 0.21
$ws.Range("F38").Value = 0.01
$ws.Range("G38").Value = -0.21
$ws.Range("H38").Value = 7.02
$ws.Range("I38").Value = 0
$ws.Range("E39").Value = 0.2
$ws.Range("F39").Value = 0.0094
$ws.Range("G39").Value = -0.2
$ws.Range("H39").Value = 8.109999999999999
$ws.Range("I39").Value = 0
$ws.Range("E40").Value = 0.25
$ws.Range("F40").Value = 0.0116
$ws.Range("G40").Value = -0.19
$ws.Range("H40").Value = 9.31
$ws.Range("I40").Value = 0
$ws.Range("E41").Value = 0.24
$ws.Range("F41").Value = 0.011
$ws.Range("G41").Value = -0.18
$ws.Range("H41").Value = 10.61
$ws.Range("E42").Value = 0.3
$ws.Range("F42").Value = 0.0136
$ws.Range("G42").Value = -0.17
$ws.Range("H42").Value = 12.02
$ws.Range("I42").Value = 0
$ws.Range("E43").Value = 0.3
$ws.Range("F43").Value = 0.0135
$ws.Range("G43").Value = -0.16
$ws.Range("H43").Value = 13.55
$ws.Range("I43").Value = 0
$ws.Range("E44").Value = 0.35
$ws.Range("F44").Value = 0.0156
$ws.Range("G44").Value = -0.16
$ws.Range("H44").Value = 15.18
$ws.Range("I44").Value = 0
$ws.Range("D45").Value = "Aguardar"
$ws.Range("E45").Value = 0.39
$ws.Range("F45").Value = 0.0171
$ws.Range("G45").Value = -0.15
$ws.Range("H45").Value = 16.92
$ws.Range("G46").Value = -0.14
$ws.Range("H46").Value = 18.76
$ws.Range("I46").Value = 0
$ws.Range("D47").Value = "Aguardar"
$ws.Range("E47").Value = 0.43
$ws.Range("F47").Value = 0.0185
$ws.Range("G47").Value = -0.13
$ws.Range("H47").Value = 20.69
$ws.Range("I47").Value = 0
$ws.Range("E48").Value = 0.51
$ws.Range("F48").Value = 0.0217
$ws.Range("G48").Value = -0.12
$ws.Range("H48").Value = 22.72
$ws.Range("I48").Value = 0
$ws.Range("D49").Value = "Montar"
$ws.Range("E49").Value = 0.5
$ws.Range("F49").Value = 0.0211
$ws.Range("G49").Value = -0.11
$ws.Range("H49").Value = 24.83
$ws.Range("I49").Value = 0
$ws.Range("E50").Value = 0.61
$ws.Range("F50").Value = 0.0254
$ws.Range("G50").Value = -0.1
$ws.Range("H50").Value = 27.02
$ws.Range("I50").Value = 0
$ws.Range("E51").Value = 0.65
$ws.Range("F51").Value = 0.0268
$ws.Range("G51").Value = -0.09
$ws.Range("H51").Value = 29.28
$ws.Range("I51").Value = 0
$ws.Range("E52").Value = 0.8
$ws.Range("F52").Value = 0.0323
$ws.Range("G52").Value = -0.07000000000000001
$ws.Range("H52").Value = 33.96
$ws.Range("I52").Value = 0
$ws.Range("E53").Value = 0.86
$ws.Range("F53").Value = 0.0344
$ws.Range("G53").Value = -0.06
$ws.Range("H53").Value = 36.37
$ws.Range("I53").Value = 0
$ws.Range("E54").Value = 0.96
$ws.Range("F54").Value = 0.038
$ws.Range("G54").Value = -0.05
$ws.Range("H54").Value = 38.81
$ws.Range("I54").Value = 0
$ws.Range("E55").Value = 1.06
$ws.Range("F55").Value = 0.0416
$ws.Range("G55").Value = -0.04
$ws.Range("H55").Value = 41.26
$ws.Range("I55").Value = 0
$ws.Range("G56").Value = -0.04
$ws.Range("H56").Value = 41.26
$ws.Range("E57").Value = 1.15
$ws.Range("F57").Value = 0.0447
$ws.Range("G57").Value = -0.03
$ws.Range("H57").Value = 43.73
$ws.Range("I57").Value = 0
$ws.Range("E58").Value = 1.23
$ws.Range("F58").Value = 0.0473
$ws.Range("G58").Value = -0.02
$ws.Range("H58").Value = 46.2
$ws.Range("I58").Value = 0
$ws.Range("E59").Value = 1.37
$ws.Range("F59").Value = 0.0522
$ws.Range("G59").Value = -0.01
$ws.Range("H59").Value = 48.66
$ws.Range("I59").Value = 0
$ws.Range("E60").Value = 1.48
$ws.Range("F60").Value = 0.0558
$ws.Range("G60").Value = -0.01
$ws.Range("H60").Value = 51.1
$ws.Range("I60").Value = 0
$ws.Range("E61").Value = 1.57
$ws.Range("F61").Value = 0.0587
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 53.51
$ws.Range("I61").Value = 0
$ws.Range("E62").Value = 1.75
$ws.Range("F62").Value = 0.0648
$ws.Range("G62").Value = 0.01
$ws.Range("H62").Value = 55.89
$ws.Range("I62").Value = 0
$ws.Range("E63").Value = 2
$ws.Range("F63").Value = 0.0727
$ws.Range("G63").Value = 0.03
$ws.Range("H63").Value = 60.51
$ws.Range("I63").Value = 0
$ws.Range("G64").Value = 0.03
$ws.Range("H64").Value = 60.51
$ws.Range("E65").Value = 2.17
$ws.Range("F65").Value = 0.07820000000000001
$ws.Range("G65").Value = 0.04
$ws.Range("H65").Value = 62.74
$ws.Range("E66").Value = 2.34
$ws.Range("F66").Value = 0.08359999999999999
$ws.Range("G66").Value = 0.05
$ws.Range("H66").Value = 64.91
$ws.Range("I66").Value = 0
$ws.Range("E67").Value = 2.58
$ws.Range("F67").Value = 0.0905
$ws.Range("G67").Value = 0.07000000000000001
$ws.Range("H67").Value = 69.05
$ws.Range("G68").Value = 0.07000000000000001
$ws.Range("H68").Value = 69.05
$ws.Range("E69").Value = 2.65
$ws.Range("F69").Value = 0.0922
$ws.Range("G69").Value = 0.08
$ws.Range("H69").Value = 71.01000000000001
$ws.Range("E70").Value = 2.9
$ws.Range("F70").Value = 0.1
$ws.Range("G70").Value = 0.09
$ws.Range("H70").Value = 72.89
$ws.Range("E71").Value = 3.28
$ws.Range("F71").Value = 0.1112
$ws.Range("G71").Value = 0.11
$ws.Range("H71").Value = 76.42
$ws.Range("E72").Value = 3.86
$ws.Range("F72").Value = 0.1287
$ws.Range("G72").Value = 0.13
$ws.Range("H72").Value = 79.63
$ws.Range("E73").Value = 4.03
$ws.Range("F73").Value = 0.1321
$ws.Range("G73").Value = 0.14
$ws.Range("H73").Value = 82.51000000000001
$ws.Range("G74").Value = 0.14
$ws.Range("H74").Value = 82.51000000000001
$ws.Range("G75").Value = 0.15
$ws.Range("H75").Value = 83.83
$ws.Range("G76").Value = 0.16
$ws.Range("H76").Value = 85.08
$ws.Range("G77").Value = 0.18
$ws.Range("H77").Value = 87.34999999999999
$ws.Range("G78").Value = 0.19
$ws.Range("H78").Value = 88.38
$ws.Range("E79").Value = 5.45
$ws.Range("F79").Value = 0.1703
$ws.Range("G79").Value = 0.2
$ws.Range("H79").Value = 89.34
$ws.Range("G80").Value = 0.22
$ws.Range("H80").Value = 91.06999999999999
$ws.Range("G81").Value = 0.24
$ws.Range("H81").Value = 92.56
